$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 14142.143  # H69: 14000 -> 14142.143
$ws.Cells.Item(69, 10).Value = 14995.5  # J69: 14996 -> 14995.5
$ws.Cells.Item(69, 12).Value = 44986.5  # L69: 44988 -> 44986.5
$ws.Cells.Item(69, 14).Value = -46734.5  # N69: -46736 -> -46734.5
$ws.Cells.Item(72, 8).Value = 14142.143  # H72: 14000 -> 14142.143
$ws.Cells.Item(72, 10).Value = 14995.5  # J72: 14996 -> 14995.5
$ws.Cells.Item(72, 12).Value = 134959.5  # L72: 134964 -> 134959.5
$ws.Cells.Item(72, 14).Value = -143695.5  # N72: -143700 -> -143695.5
$ws.Cells.Item(112, 8).Value = 10627.25  # H112: 12329.8 -> 10627.25
$ws.Cells.Item(112, 9).Value = 1124.75  # I112: 1159.8 -> 1124.75
$ws.Cells.Item(112, 10).Value = 15378.5  # J112: 23499.8 -> 15378.5
$ws.Cells.Item(112, 11).Value = 3374.25  # K112: 3479.4 -> 3374.25
$ws.Cells.Item(112, 12).Value = 46135.5  # L112: 70499.39999999999 -> 46135.5
$ws.Cells.Item(112, 13).Value = -2266.25  # M112: -2371.4 -> -2266.25
$ws.Cells.Item(112, 14).Value = -48351.5  # N112: -72715.39999999999 -> -48351.5
$ws.Cells.Item(132, 8).Value = 1963.8451  # H132: 1972.8472 -> 1963.8451
$ws.Cells.Item(132, 9).Value = 1949.0834  # I132: 1959.9508 -> 1949.0834
$ws.Cells.Item(132, 11).Value = 5847.2502  # K132: 5879.8524 -> 5847.2502
$ws.Cells.Item(132, 13).Value = -3317.2502  # M132: -3349.8524 -> -3317.2502
$ws.Cells.Item(137, 8).Value = 1100.1428  # H137: 1131 -> 1100.1428
$ws.Cells.Item(137, 9).Value = 887.5454999999999  # I137: 906.4 -> 887.5454999999999
$ws.Cells.Item(137, 11).Value = 2662.6365  # K137: 2719.2 -> 2662.6365
$ws.Cells.Item(137, 13).Value = -112.6364999999996  # M137: -169.1999999999998 -> -112.6364999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 2928.75  # H4: 2870.3333 -> 2928.75
$ws.Cells.Item(4, 9).Value = 1481.125  # I4: 1383.1111 -> 1481.125
$ws.Cells.Item(4, 10).Value = 5824  # J4: 7332 -> 5824
$ws.Cells.Item(4, 11).Value = 1481.125  # K4: 1383.1111 -> 1481.125
$ws.Cells.Item(4, 12).Value = 5824  # L4: 7332 -> 5824
$ws.Cells.Item(4, 13).Value = -1365.125  # M4: -1267.1111 -> -1365.125
$ws.Cells.Item(4, 14).Value = -6056  # N4: -7564 -> -6056
$ws.Cells.Item(6, 8).Value = 10001835  # H6: 10002499 -> 10001835
$ws.Cells.Item(6, 9).Value = 2002  # I6: 0 -> 2002
$ws.Cells.Item(6, 10).Value = 30001500  # J6: 10002499 -> 30001500
$ws.Cells.Item(6, 11).Value = 2002  # K6: 0 -> 2002
$ws.Cells.Item(6, 12).Value = 30001500  # L6: 10002499 -> 30001500
$ws.Cells.Item(6, 13).Value = -1829  # M6: None -> -1829
$ws.Cells.Item(6, 14).Value = -30001846  # N6: -10002845 -> -30001846
$ws.Cells.Item(23, 8).Value = 20000  # H23: 0 -> 20000
$ws.Cells.Item(23, 9).Value = 20000  # I23: 0 -> 20000
$ws.Cells.Item(23, 11).Value = 20000  # K23: 0 -> 20000
$ws.Cells.Item(23, 13).Value = -19741  # M23: None -> -19741
$ws.Cells.Item(37, 8).Value = 29499.5  # H37: 16311.667 -> 29499.5
$ws.Cells.Item(37, 10).Value = 46999  # J37: 24935 -> 46999
$ws.Cells.Item(37, 12).Value = 46999  # L37: 24935 -> 46999
$ws.Cells.Item(37, 14).Value = -47545  # N37: -25481 -> -47545
$ws.Cells.Item(63, 8).Value = 500  # H63: 999.1111 -> 500
$ws.Cells.Item(63, 9).Value = 500  # I63: 873.25 -> 500
$ws.Cells.Item(63, 10).Value = 0  # J63: 2006 -> 0
$ws.Cells.Item(63, 11).Value = 500  # K63: 873.25 -> 500
$ws.Cells.Item(63, 12).Value = 0  # L63: 2006 -> 0
$ws.Cells.Item(63, 13).Value = $null  # M63: -187.25 -> None
$ws.Cells.Item(63, 14).Value = 186  # N63: -3378 -> 186
$ws.Cells.Item(66, 8).Value = 500  # H66: 999.1111 -> 500
$ws.Cells.Item(66, 9).Value = 500  # I66: 873.25 -> 500
$ws.Cells.Item(66, 10).Value = 0  # J66: 2006 -> 0
$ws.Cells.Item(66, 11).Value = 2500  # K66: 4366.25 -> 2500
$ws.Cells.Item(66, 12).Value = 0  # L66: 10030 -> 0
$ws.Cells.Item(66, 13).Value = $null  # M66: -934.25 -> None
$ws.Cells.Item(66, 14).Value = 932  # N66: -16894 -> 932
$ws.Cells.Item(74, 8).Value = 3000  # H74: 0 -> 3000
$ws.Cells.Item(74, 10).Value = 3000  # J74: 0 -> 3000
$ws.Cells.Item(74, 12).Value = 3000  # L74: 0 -> 3000
$ws.Cells.Item(74, 14).Value = -4748  # N74: None -> -4748
$ws.Cells.Item(77, 8).Value = 3000  # H77: 0 -> 3000
$ws.Cells.Item(77, 10).Value = 3000  # J77: 0 -> 3000
$ws.Cells.Item(77, 12).Value = 15000  # L77: 0 -> 15000
$ws.Cells.Item(77, 14).Value = -23736  # N77: None -> -23736
$ws.Cells.Item(97, 8).Value = 2332.7693  # H97: 2212.5715 -> 2332.7693
$ws.Cells.Item(97, 9).Value = 932.6  # I97: 906.9091 -> 932.6
$ws.Cells.Item(97, 11).Value = 932.6  # K97: 906.9091 -> 932.6
$ws.Cells.Item(97, 13).Value = -436.6  # M97: -410.9091 -> -436.6
$ws.Cells.Item(135, 8).Value = 70979.8  # H135: 72299.8 -> 70979.8
$ws.Cells.Item(135, 10).Value = 70979.8  # J135: 72299.8 -> 70979.8
$ws.Cells.Item(135, 12).Value = 70979.8  # L135: 72299.8 -> 70979.8
$ws.Cells.Item(135, 14).Value = -81119.8  # N135: -82439.8 -> -81119.8
$ws.Cells.Item(139, 8).Value = 85239.7  # H139: 85247.2 -> 85239.7
$ws.Cells.Item(139, 10).Value = 85239.7  # J139: 85247.2 -> 85239.7
$ws.Cells.Item(139, 12).Value = 85239.7  # L139: 85247.2 -> 85239.7
$ws.Cells.Item(139, 14).Value = -95519.7  # N139: -95527.2 -> -95519.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 0  # H15: 5900 -> 0
$ws.Cells.Item(15, 10).Value = 0  # J15: 5900 -> 0
$ws.Cells.Item(15, 12).Value = $null  # L15: 5900 -> None
$ws.Cells.Item(15, 14).Value = 0  # N15: -6354 -> 0
$ws.Cells.Item(20, 8).Value = 2178.7273  # H20: 2097.0833 -> 2178.7273
$ws.Cells.Item(20, 9).Value = 1974.625  # I20: 1888.4445 -> 1974.625
$ws.Cells.Item(20, 11).Value = 1974.625  # K20: 1888.4445 -> 1974.625
$ws.Cells.Item(20, 13).Value = -1727.625  # M20: -1641.4445 -> -1727.625
$ws.Cells.Item(35, 8).Value = 60585  # H35: 60585.332 -> 60585
$ws.Cells.Item(35, 10).Value = 65033.125  # J35: 65033.5 -> 65033.125
$ws.Cells.Item(35, 12).Value = 65033.125  # L35: 65033.5 -> 65033.125
$ws.Cells.Item(35, 14).Value = -65653.125  # N35: -65653.5 -> -65653.125
$ws.Cells.Item(81, 8).Value = 25893.2  # H81: 20666.285 -> 25893.2
$ws.Cells.Item(81, 10).Value = 25893.2  # J81: 20666.285 -> 25893.2
$ws.Cells.Item(81, 12).Value = 25893.2  # L81: 20666.285 -> 25893.2
$ws.Cells.Item(81, 14).Value = -28015.2  # N81: -22788.285 -> -28015.2
$ws.Cells.Item(84, 8).Value = 25893.2  # H84: 20666.285 -> 25893.2
$ws.Cells.Item(84, 10).Value = 25893.2  # J84: 20666.285 -> 25893.2
$ws.Cells.Item(84, 12).Value = 77679.60000000001  # L84: 61998.855 -> 77679.60000000001
$ws.Cells.Item(84, 14).Value = -88287.60000000001  # N84: -72606.855 -> -88287.60000000001
$ws.Cells.Item(105, 8).Value = 45458196  # H105: 25002590 -> 45458196
$ws.Cells.Item(105, 9).Value = 71429720  # I105: 33334552 -> 71429720
$ws.Cells.Item(105, 10).Value = 8027.75  # J105: 6701.6 -> 8027.75
$ws.Cells.Item(105, 11).Value = 71429720  # K105: 33334552 -> 71429720
$ws.Cells.Item(105, 12).Value = 8027.75  # L105: 6701.6 -> 8027.75
$ws.Cells.Item(105, 13).Value = -71427973  # M105: -33332805 -> -71427973
$ws.Cells.Item(105, 14).Value = -11521.75  # N105: -10195.6 -> -11521.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 56198.6  # H59: 63498.668 -> 56198.6
$ws.Cells.Item(59, 10).Value = 53997.668  # J59: 65498 -> 53997.668
$ws.Cells.Item(59, 12).Value = 53997.668  # L59: 65498 -> 53997.668
$ws.Cells.Item(59, 14).Value = -56287.668  # N59: -67788 -> -56287.668
$ws.Cells.Item(68, 8).Value = 57313.418  # H68: 56978.46 -> 57313.418
$ws.Cells.Item(68, 10).Value = 57313.418  # J68: 56978.46 -> 57313.418
$ws.Cells.Item(68, 12).Value = 57313.418  # L68: 56978.46 -> 57313.418
$ws.Cells.Item(68, 14).Value = -58811.418  # N68: -58476.46 -> -58811.418
$ws.Cells.Item(71, 8).Value = 57313.418  # H71: 56978.46 -> 57313.418
$ws.Cells.Item(71, 10).Value = 57313.418  # J71: 56978.46 -> 57313.418
$ws.Cells.Item(71, 12).Value = 171940.254  # L71: 170935.38 -> 171940.254
$ws.Cells.Item(71, 14).Value = -179428.254  # N71: -178423.38 -> -179428.254
$ws.Cells.Item(107, 8).Value = 1358.4517  # H107: 1383.2333 -> 1358.4517
$ws.Cells.Item(107, 9).Value = 774.12  # I107: 780.75 -> 774.12
$ws.Cells.Item(107, 11).Value = 774.12  # K107: 780.75 -> 774.12
$ws.Cells.Item(107, 13).Value = 1145.88  # M107: 1139.25 -> 1145.88

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1245.4166  # H68: 1268.7273 -> 1245.4166
$ws.Cells.Item(68, 10).Value = 1643.25  # J68: 1861.3334 -> 1643.25
$ws.Cells.Item(68, 12).Value = 4929.75  # L68: 5584.0002 -> 4929.75
$ws.Cells.Item(68, 14).Value = -6551.75  # N68: -7206.0002 -> -6551.75
$ws.Cells.Item(71, 8).Value = 1245.4166  # H71: 1268.7273 -> 1245.4166
$ws.Cells.Item(71, 10).Value = 1643.25  # J71: 1861.3334 -> 1643.25
$ws.Cells.Item(71, 12).Value = 14789.25  # L71: 16752.0006 -> 14789.25
$ws.Cells.Item(71, 14).Value = -22901.25  # N71: -24864.0006 -> -22901.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7094.533  # H70: 7245.0713 -> 7094.533
$ws.Cells.Item(70, 9).Value = 6545.8887  # I70: 6740.75 -> 6545.8887
$ws.Cells.Item(70, 11).Value = 6545.8887  # K70: 6740.75 -> 6545.8887
$ws.Cells.Item(70, 13).Value = -6275.8887  # M70: -6470.75 -> -6275.8887
$ws.Cells.Item(73, 8).Value = 7094.533  # H73: 7245.0713 -> 7094.533
$ws.Cells.Item(73, 9).Value = 6545.8887  # I73: 6740.75 -> 6545.8887
$ws.Cells.Item(73, 11).Value = 6545.8887  # K73: 6740.75 -> 6545.8887
$ws.Cells.Item(73, 13).Value = -5609.8887  # M73: -5804.75 -> -5609.8887
$ws.Cells.Item(122, 8).Value = 80391.92  # H122: 53240.25 -> 80391.92
$ws.Cells.Item(122, 9).Value = 103266.3  # I122: 61904.293 -> 103266.3
$ws.Cells.Item(122, 11).Value = 309798.9  # K122: 185712.879 -> 309798.9
$ws.Cells.Item(122, 13).Value = -307348.9  # M122: -183262.879 -> -307348.9

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 36451.285  # H50: 30124.584 -> 36451.285
$ws.Cells.Item(50, 10).Value = 43599.8  # J50: 32433.5 -> 43599.8
$ws.Cells.Item(50, 12).Value = 43599.8  # L50: 32433.5 -> 43599.8
$ws.Cells.Item(50, 14).Value = -44873.8  # N50: -33707.5 -> -44873.8
$ws.Cells.Item(54, 8).Value = 0  # H54: 39084 -> 0
$ws.Cells.Item(54, 10).Value = 0  # J54: 39084 -> 0
$ws.Cells.Item(54, 12).Value = $null  # L54: 39084 -> None
$ws.Cells.Item(54, 14).Value = 0  # N54: -40372 -> 0
$ws.Cells.Item(56, 8).Value = 30129.6  # H56: 36316 -> 30129.6
$ws.Cells.Item(56, 9).Value = 30129.6  # I56: 36316 -> 30129.6
$ws.Cells.Item(56, 11).Value = 30129.6  # K56: 36316 -> 30129.6
$ws.Cells.Item(56, 13).Value = -29438.6  # M56: -35625 -> -29438.6
$ws.Cells.Item(68, 8).Value = 2672.3635  # H68: 2672.2727 -> 2672.3635
$ws.Cells.Item(68, 9).Value = 2174.625  # I68: 2174.5 -> 2174.625
$ws.Cells.Item(68, 11).Value = 2174.625  # K68: 2174.5 -> 2174.625
$ws.Cells.Item(68, 13).Value = -1425.625  # M68: -1425.5 -> -1425.625
$ws.Cells.Item(71, 8).Value = 2672.3635  # H71: 2672.2727 -> 2672.3635
$ws.Cells.Item(71, 9).Value = 2174.625  # I71: 2174.5 -> 2174.625
$ws.Cells.Item(71, 11).Value = 10873.125  # K71: 10872.5 -> 10873.125
$ws.Cells.Item(71, 13).Value = -7129.125  # M71: -7128.5 -> -7129.125
$ws.Cells.Item(93, 8).Value = 1974.4  # H93: 1624.8462 -> 1974.4
$ws.Cells.Item(93, 9).Value = 2104.889  # I93: 1693.5834 -> 2104.889
$ws.Cells.Item(93, 11).Value = 2104.889  # K93: 1693.5834 -> 2104.889
$ws.Cells.Item(93, 13).Value = -856.8890000000001  # M93: -445.5834 -> -856.8890000000001
$ws.Cells.Item(122, 8).Value = 3614.3333  # H122: 3517.2415 -> 3614.3333
$ws.Cells.Item(122, 9).Value = 3381.125  # I122: 3347.739 -> 3381.125
$ws.Cells.Item(122, 10).Value = 4360.6  # J122: 4167 -> 4360.6
$ws.Cells.Item(122, 11).Value = 10143.375  # K122: 10043.217 -> 10143.375
$ws.Cells.Item(122, 12).Value = 13081.8  # L122: 12501 -> 13081.8
$ws.Cells.Item(122, 13).Value = -7693.375  # M122: -7593.217000000001 -> -7693.375
$ws.Cells.Item(122, 14).Value = -17981.8  # N122: -17401 -> -17981.8
$ws.Cells.Item(136, 8).Value = 9527876  # H136: 8337402 -> 9527876
$ws.Cells.Item(136, 9).Value = 11908113  # I136: 12349066 -> 11908113
$ws.Cells.Item(136, 10).Value = 6928.5713  # J136: 5484.846 -> 6928.5713
$ws.Cells.Item(136, 11).Value = 35724339  # K136: 37047198 -> 35724339
$ws.Cells.Item(136, 12).Value = 20785.7139  # L136: 16454.538 -> 20785.7139
$ws.Cells.Item(136, 13).Value = -35721789  # M136: -37044648 -> -35721789
$ws.Cells.Item(136, 14).Value = -25885.7139  # N136: -21554.538 -> -25885.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(61, 8).Value = 26595  # H61: 27081.428 -> 26595
$ws.Cells.Item(61, 10).Value = 37792.5  # J61: 36234 -> 37792.5
$ws.Cells.Item(61, 12).Value = 37792.5  # L61: 36234 -> 37792.5
$ws.Cells.Item(61, 14).Value = -38376.5  # N61: -36818 -> -38376.5
$ws.Cells.Item(62, 8).Value = 17358.666  # H62: 17158.533 -> 17358.666
$ws.Cells.Item(62, 9).Value = 16438.3  # I62: 16398.092 -> 16438.3
$ws.Cells.Item(62, 10).Value = 19199.4  # J62: 19249.75 -> 19199.4
$ws.Cells.Item(62, 11).Value = 16438.3  # K62: 16398.092 -> 16438.3
$ws.Cells.Item(62, 12).Value = 19199.4  # L62: 19249.75 -> 19199.4
$ws.Cells.Item(62, 13).Value = -15814.3  # M62: -15774.092 -> -15814.3
$ws.Cells.Item(62, 14).Value = -20447.4  # N62: -20497.75 -> -20447.4
$ws.Cells.Item(65, 8).Value = 17358.666  # H65: 17158.533 -> 17358.666
$ws.Cells.Item(65, 9).Value = 16438.3  # I65: 16398.092 -> 16438.3
$ws.Cells.Item(65, 10).Value = 19199.4  # J65: 19249.75 -> 19199.4
$ws.Cells.Item(65, 11).Value = 82191.5  # K65: 81990.46000000001 -> 82191.5
$ws.Cells.Item(65, 12).Value = 95997  # L65: 96248.75 -> 95997
$ws.Cells.Item(65, 13).Value = -79071.5  # M65: -78870.46000000001 -> -79071.5
$ws.Cells.Item(65, 14).Value = -102237  # N65: -102488.75 -> -102237
$ws.Cells.Item(70, 8).Value = 54000  # H70: 53368.332 -> 54000
$ws.Cells.Item(70, 10).Value = 54000  # J70: 53368.332 -> 54000
$ws.Cells.Item(70, 12).Value = 54000  # L70: 53368.332 -> 54000
$ws.Cells.Item(70, 14).Value = -54630  # N70: -53998.332 -> -54630
$ws.Cells.Item(73, 8).Value = 54000  # H73: 53368.332 -> 54000
$ws.Cells.Item(73, 10).Value = 54000  # J73: 53368.332 -> 54000
$ws.Cells.Item(73, 12).Value = 54000  # L73: 53368.332 -> 54000
$ws.Cells.Item(73, 14).Value = -56184  # N73: -55552.332 -> -56184
$ws.Cells.Item(96, 8).Value = 2033.05  # H96: 2071.6316 -> 2033.05
$ws.Cells.Item(96, 9).Value = 2430  # I96: 2712.5 -> 2430
$ws.Cells.Item(96, 11).Value = 2430  # K96: 2712.5 -> 2430
$ws.Cells.Item(96, 13).Value = -1057  # M96: -1339.5 -> -1057
